$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 25.00841100000001
$ws.Range("H2").Value = 75.02523300000001
$ws.Range("I2").Value = 0.4156829172908309
$ws.Range("J2").Value = 0.415682917290831
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.457778999999999
$ws.Range("N2").Value = 22.373337
$ws.Range("O2").Value = 0.08029647035915141
$ws.Range("P2").Value = 0.0802964703591514
$ws.Range("Q2").Value = 186.507202379169
$ws.Range("R2").Value = 1678.564821412521
$ws.Range("S2").Value = 0.03337787104704879
$ws.Range("T2").Value = 0.03337787104704879

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 25.00841100000001
$ws.Range("H3").Value = 75.02523300000001
$ws.Range("I3").Value = 0.4156829172908309
$ws.Range("J3").Value = 0.415682917290831
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.732509
$ws.Range("N3").Value = 8.197527
$ws.Range("O3").Value = 0.02942039820764526
$ws.Range("P3").Value = 0.02942039820764525
$ws.Range("Q3").Value = 68.335708133199
$ws.Range("R3").Value = 615.021373198791
$ws.Range("S3").Value = 0.01222955695481192
$ws.Range("T3").Value = 0.01222955695481191

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 25.00841100000001
$ws.Range("H4").Value = 75.02523300000001
$ws.Range("I4").Value = 0.4156829172908309
$ws.Range("J4").Value = 0.415682917290831
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 82.68775466666666
$ws.Range("N4").Value = 248.063264
$ws.Range("O4").Value = 0.8902831314332034
$ws.Range("P4").Value = 0.8902831314332033
$ws.Range("Q4").Value = 2067.889353371168
$ws.Range("R4").Value = 18611.00418034051
$ws.Range("S4").Value = 0.3700754892889702
$ws.Range("T4").Value = 0.3700754892889702

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 20.496019
$ws.Range("H5").Value = 61.488057
$ws.Range("I5").Value = 0.340679180727168
$ws.Range("J5").Value = 0.3406791807271681
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.457778999999999
$ws.Range("N5").Value = 22.373337
$ws.Range("O5").Value = 0.08029647035915141
$ws.Range("P5").Value = 0.0802964703591514
$ws.Range("Q5").Value = 152.854780081801
$ws.Range("R5").Value = 1375.693020736209
$ws.Range("S5").Value = 0.02735533573723904
$ws.Range("T5").Value = 0.02735533573723903

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 20.496019
$ws.Range("H6").Value = 61.488057
$ws.Range("I6").Value = 0.340679180727168
$ws.Range("J6").Value = 0.3406791807271681
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.732509
$ws.Range("N6").Value = 8.197527
$ws.Range("O6").Value = 0.02942039820764526
$ws.Range("P6").Value = 0.02942039820764525
$ws.Range("Q6").Value = 56.005556381671
$ws.Range("R6").Value = 504.0500074350389
$ws.Range("S6").Value = 0.01002291715804763
$ws.Range("T6").Value = 0.01002291715804763

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 20.496019
$ws.Range("H7").Value = 61.488057
$ws.Range("I7").Value = 0.340679180727168
$ws.Range("J7").Value = 0.3406791807271681
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 82.68775466666666
$ws.Range("N7").Value = 248.063264
$ws.Range("O7").Value = 0.8902831314332034
$ws.Range("P7").Value = 0.8902831314332033
$ws.Range("Q7").Value = 1694.769790715339
$ws.Range("R7").Value = 15252.92811643805
$ws.Range("S7").Value = 0.3033009278318814
$ws.Range("T7").Value = 0.3033009278318814

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.65779933333333
$ws.Range("H8").Value = 43.973398
$ws.Range("I8").Value = 0.243637901982001
$ws.Range("J8").Value = 0.243637901982001
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.457778999999999
$ws.Range("N8").Value = 22.373337
$ws.Range("O8").Value = 0.08029647035915141
$ws.Range("P8").Value = 0.0802964703591514
$ws.Range("Q8").Value = 109.3146280543473
$ws.Range("R8").Value = 983.831652489126
$ws.Range("S8").Value = 0.01956326357486358
$ws.Range("T8").Value = 0.01956326357486358

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.65779933333333
$ws.Range("H9").Value = 43.973398
$ws.Range("I9").Value = 0.243637901982001
$ws.Range("J9").Value = 0.243637901982001
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.732509
$ws.Range("N9").Value = 8.197527
$ws.Range("O9").Value = 0.02942039820764526
$ws.Range("P9").Value = 0.02942039820764525
$ws.Range("Q9").Value = 40.05256859852734
$ws.Range("R9").Value = 360.473117386746
$ws.Range("S9").Value = 0.007167924094785714
$ws.Range("T9").Value = 0.007167924094785712

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.65779933333333
$ws.Range("H10").Value = 43.973398
$ws.Range("I10").Value = 0.243637901982001
$ws.Range("J10").Value = 0.243637901982001
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 82.68775466666666
$ws.Range("N10").Value = 248.063264
$ws.Range("O10").Value = 0.8902831314332034
$ws.Range("P10").Value = 0.8902831314332033
$ws.Range("Q10").Value = 1212.020515227897
$ws.Range("R10").Value = 10908.18463705107
$ws.Range("S10").Value = 0.2169067143123517
$ws.Range("T10").Value = 0.2169067143123517
